$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values on row 102 (2025-02-26, abs_activity)
$ws.Range("C102").Value = 8.172556287543546
$ws.Range("F102").Value = 8.172556287543546

# Append new rows for 2025-02-27 (abs_activity, rel_activity, abs_sleep, rel_sleep)
# Force column A to a text format before assignment so the date-like
# strings are not auto-converted into date serial numbers, then restore
# the default "Normal" style so no residual cell formatting is left behind.
$ws.Range("A106:A109").NumberFormat = "@"

$ws.Range("A106").Value = "2025-02-27"
$ws.Range("B106").Value = "abs_activity"
$ws.Range("C106").Value = 9.614946288626287
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 0
$ws.Range("F106").Value = 9.614946288626287

$ws.Range("A107").Value = "2025-02-27"
$ws.Range("B107").Value = "rel_activity"
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 0
$ws.Range("E107").Value = 0
$ws.Range("F107").Value = 0

$ws.Range("A108").Value = "2025-02-27"
$ws.Range("B108").Value = "abs_sleep"
$ws.Range("C108").Value = 7.733333333333334
$ws.Range("D108").Value = 0
$ws.Range("E108").Value = 0
$ws.Range("F108").Value = 7.733333333333334

$ws.Range("A109").Value = "2025-02-27"
$ws.Range("B109").Value = "rel_sleep"
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 0
$ws.Range("E109").Value = 0
$ws.Range("F109").Value = 0

$ws.Range("A106:A109").Style = "Normal"
